# Auto-generated edit script applying latest crypto price/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '79.702.65'
$ws.Range("E2").Value = '  +4.36%  '

$ws.Range("D3").Value = '3.206.10'
$ws.Range("E3").Value = '  +5.31%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '636.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.06%  '

$ws.Range("E8").Value = '  +19.95%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.603'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.90%  '

$ws.Range("D10").Value = '3.204.65'
$ws.Range("E10").Value = '  +5.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.595'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +35.74%  '

$ws.Range("E12").Value = '  +35.43%  '

$ws.Range("E13").Value = '  +3.16%  '

$ws.Range("E14").Value = '  +3.06%  '

$ws.Range("D15").Value = '3.793.41'
$ws.Range("E15").Value = '  +5.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +12.06%  '

$ws.Range("D17").Value = '79.519.21'
$ws.Range("E17").Value = '  +4.21%  '

$ws.Range("D18").Value = '3.201.92'
$ws.Range("E18").Value = '  +5.03%  '

$ws.Range("E19").Value = '  +7.30%  '

$ws.Range("E20").Value = '  +25.73%  '

$ws.Range("E21").Value = '  +4.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '440.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +15.28%  '

$ws.Range("E23").Value = '  +19.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.59%  '

$ws.Range("D25").Value = '3.368.24'
$ws.Range("E25").Value = '  +5.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '77.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.42%  '

$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000124'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +15.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.39%  '

$ws.Range("E32").Value = '  +9.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '560.62'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +13.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.157'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +33.96%  '

$ws.Range("E35").Value = '  +6.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.121'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +16.21%  '

$ws.Range("E38").Value = '  -0.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.414'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '163.38'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '20.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.96%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '192.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("E45").Value = '  +11.85%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.801'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '43.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.55%  '

$ws.Range("E50").Value = '  +7.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.65'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +15.56%  '
